# Fruta / hortaliza, semanal
# Insert two new weekly price observations into the "Apio" data block.
# The sheet's rows are not sorted by date, so the two newly observed
# days (2022-01-06 and 2022-01-07) land in the middle of the existing
# block, pushing everything below each insertion point down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First insertion: new row at (original) row 188 -------------------
# Shifts old rows 188..304 down to 189..305.
$ws.Rows("188:188").Insert()

$ws.Cells.Item(188, 1).Value  = 3
$ws.Cells.Item(188, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(188, 3).Value  = "Coquimbo"
$ws.Cells.Item(188, 4).Value  = 44567
$ws.Cells.Item(188, 5).Value  = 5
$ws.Cells.Item(188, 6).Value  = 100112017
$ws.Cells.Item(188, 7).Value  = "Apio"
$ws.Cells.Item(188, 8).Value  = "Americana (o)"
$ws.Cells.Item(188, 9).Value  = "Primera"
$ws.Cells.Item(188, 10).Value = 160
$ws.Cells.Item(188, 11).Value = 9000
$ws.Cells.Item(188, 12).Value = 9000
$ws.Cells.Item(188, 13).Value = 9000
$ws.Cells.Item(188, 14).Value = "`$/docena de matas"
$ws.Cells.Item(188, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(188, 16).Value = 1500
$ws.Cells.Item(188, 17).Value = 6
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# --- Second insertion: new row at row 294 (post first-insert numbering)
# Shifts rows 294..305 down to 295..306.
$ws.Rows("294:294").Insert()

$ws.Cells.Item(294, 1).Value  = 3
$ws.Cells.Item(294, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(294, 3).Value  = "Coquimbo"
$ws.Cells.Item(294, 4).Value  = 44568
$ws.Cells.Item(294, 5).Value  = 5
$ws.Cells.Item(294, 6).Value  = 100112017
$ws.Cells.Item(294, 7).Value  = "Apio"
$ws.Cells.Item(294, 8).Value  = "Americana (o)"
$ws.Cells.Item(294, 9).Value  = "Primera"
$ws.Cells.Item(294, 10).Value = 160
$ws.Cells.Item(294, 11).Value = 9000
$ws.Cells.Item(294, 12).Value = 9000
$ws.Cells.Item(294, 13).Value = 9000
$ws.Cells.Item(294, 14).Value = "`$/docena de matas"
$ws.Cells.Item(294, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(294, 16).Value = 1500
$ws.Cells.Item(294, 17).Value = 6
$ws.Cells.Item(294, 18).Value = "Hortaliza"
